# ---------------------------------------------------------------------------
# Szkolenie warsztatowe Dapper - add the "Cwiczenia" (exercises) section
# (1 overview slide + 5 exercise slides) right before the closing slide, and
# refresh the cached "date updated automatically" placeholder text that ships
# on the slide master / all slide layouts.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation
$contentLayout = 2   # "Tytul i zawartosc" (Title and Content) layout

# --- Ćwiczenia ---
$s1 = $p.Slides.Add($p.Slides.Count, $contentLayout)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Wyszukanie klasy ćwiczenia (Exercise_X, gdzie X jest numerem ćwiczenia)`rW metodzie „RunExercise” (oraz w konstruktorze jeśli istnieję), należy zaimplementować wskazane ćwiczenie`rSprawdzenie wykonania ćwiczenia (RunUnitTests na projekcie ExerciseTests bądź na pojedynczej klasie testu)"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Ćwiczenia"

# --- Ćwiczenie 1 ---
$s2 = $p.Slides.Add($p.Slides.Count, $contentLayout)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Wykonanie zapytania i zwrócenie JEDNEGO obiektu typu EmployeeDTO z bazy danych`rTabela „Address”`r"
$s2.Shapes.Item(1).TextFrame.TextRange.Characters(33, 7).Font.Bold = $true
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "Ćwiczenie 1"

# --- Ćwiczenie 2 ---
$s3 = $p.Slides.Add($p.Slides.Count, $contentLayout)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Wykonanie procedury AddressSave`rPrzekazanie anonimowych parametrów`rWyjście: liczba dodanych adresów"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "Ćwiczenie 2"

# --- Ćwiczenie 3 ---
$s4 = $p.Slides.Add($p.Slides.Count, $contentLayout)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Wykonanie procedury AddressUpdate`rPrzekazanie parametru obiektowego (AddressDTO)`rWyjście: liczba zaktualizowanych adresów"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "Ćwiczenie 3"

# --- Ćwiczenie 4 ---
$s5 = $p.Slides.Add($p.Slides.Count, $contentLayout)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Konstruktor: Inicjalizacja listy/tablicy liczb (np. Enumerable.Range(1,2))`rMetoda: Zapytanie z tabeli adresów z przekazaniem parametru w postaci wcześniej utworzonej listy/tablicy liczb`rWyjście: Lista AddressDTO"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "Ćwiczenie 4"

# --- Ćwiczenie 5 ---
$s6 = $p.Slides.Add($p.Slides.Count, $contentLayout)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Utworzenie multi-mapowania jeden do jednego, tabeli pracowników z użytkownikami`rZapytanie powinno „joinować” tabelę pracowników z użytkownikami`rWyjście: EmployeeDTO z uzupełnionym obiektem UserDTO"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Ćwiczenie 5"

# ---------------------------------------------------------------------------
# Refresh the cached date text ("12.11.2017" -> "20.11.2017") wherever the
# datetimeFigureOut placeholder appears: the slide master and every layout.
# ---------------------------------------------------------------------------
function Update-CachedDate($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "12.11.2017") {
                $tr.Text = "20.11.2017"
            }
        }
    }
}

$master = $p.SlideMaster
Update-CachedDate $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-CachedDate $master.CustomLayouts.Item($li).Shapes
}
